$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 165-166, pushing the existing rows 165-207 down to 167-209.
$ws.Range("A165:A166").EntireRow.Insert()

# Row 165: new weekly "Primera" entry dated 2022-07-12 (serial 44754),
# same market/quality/price data as the (now shifted) neighbouring rows.
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = "2022-07-12"
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = 100112040
$ws.Cells.Item(165, 7).Value = "Cilantro"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 200
$ws.Cells.Item(165, 11).Value = 600
$ws.Cells.Item(165, 12).Value = 700
$ws.Cells.Item(165, 13).Value = 650
$ws.Cells.Item(165, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(165, 15).Value = "Región de Ñuble"
$ws.Cells.Item(165, 16).Value = 650
$ws.Cells.Item(165, 17).Value = 1
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# Row 166: matching "Segunda" entry, same date.
$ws.Cells.Item(166, 1).Value = 11
$ws.Cells.Item(166, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(166, 3).Value = "Bíobío"
$ws.Cells.Item(166, 4).Value = "2022-07-12"
$ws.Cells.Item(166, 5).Value = 8
$ws.Cells.Item(166, 6).Value = 100112040
$ws.Cells.Item(166, 7).Value = "Cilantro"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Segunda"
$ws.Cells.Item(166, 10).Value = 100
$ws.Cells.Item(166, 11).Value = 500
$ws.Cells.Item(166, 12).Value = 500
$ws.Cells.Item(166, 13).Value = 500
$ws.Cells.Item(166, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(166, 15).Value = "Región de Ñuble"
$ws.Cells.Item(166, 16).Value = 500
$ws.Cells.Item(166, 17).Value = 1
$ws.Cells.Item(166, 18).Value = "Hortaliza"

Write-Host "Final used range:" $ws.UsedRange.Address()
